# Adds a "TEST VIBRATION" row (label key + English caption) to the
# "uso comune" translations sheet, right before the existing "SUNDAY" row
# (new row 58), shifting the rows below it down by one.
# Also restores the natural selection/active-tab state that Excel leaves
# behind after this editing session:
#   - "Misc" sheet: selection left on A32 (no longer whole-sheet selected)
#   - "Selections" sheet: no longer the active tab / no longer scrolled
#   - "uso comune" sheet: becomes the active tab, selection on the new
#     STOP row (A59, after the insert)

$wb = $excel.ActiveWorkbook

# --- touch the "Misc" sheet's selection first (without leaving it active) ---
$wsMisc = $wb.Worksheets.Item("Misc")
$wsMisc.Range("A32").Select()

# --- main edit: insert the new translation row in "uso comune" ---
$wsUso = $wb.Worksheets.Item("uso comune")
$wsUso.Activate()

$wsUso.Rows("58:58").Insert()
$wsUso.Rows("58:58").RowHeight = 20.1

$wsUso.Range("A58").Value = '$LAB_STOP'
$wsUso.Range("B58").Value = 'STOP'

# leave the selection on the row that used to be row 58 (now row 59)
$wsUso.Range("A59").Select()
